$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows continue the existing daily log (columns: date, weekday, hour, rank).
# Column A holds date-like text (e.g. "2026/12/29"), matching the rest of the
# sheet where dates are stored as plain strings, not Excel date serials. Excel
# auto-converts such text to a date on assignment, so pre-format the whole new
# block as Text ("@") before writing, then restore the default "Normal" style
# (the source rows carry no explicit style) once the values are in place.
$dateRange = $ws.Range("A3044:A3085")
$dateRange.NumberFormat = "@"

$ws.Cells.Item(3044, 1).Value = "2026/12/29"
$ws.Cells.Item(3044, 2).Value = "火"
$ws.Cells.Item(3044, 3).Value = 13
$ws.Cells.Item(3044, 4).Value = 201
$ws.Cells.Item(3045, 1).Value = "2026/12/29"
$ws.Cells.Item(3045, 2).Value = "火"
$ws.Cells.Item(3045, 3).Value = 16
$ws.Cells.Item(3045, 4).Value = 201
$ws.Cells.Item(3046, 1).Value = "2026/12/29"
$ws.Cells.Item(3046, 2).Value = "火"
$ws.Cells.Item(3046, 3).Value = 19
$ws.Cells.Item(3046, 4).Value = 201
$ws.Cells.Item(3047, 1).Value = "2026/12/29"
$ws.Cells.Item(3047, 2).Value = "火"
$ws.Cells.Item(3047, 3).Value = 23
$ws.Cells.Item(3047, 4).Value = 201
$ws.Cells.Item(3048, 1).Value = "2026/12/30"
$ws.Cells.Item(3048, 2).Value = "水"
$ws.Cells.Item(3048, 3).Value = 2
$ws.Cells.Item(3048, 4).Value = 201
$ws.Cells.Item(3049, 1).Value = "2026/12/30"
$ws.Cells.Item(3049, 2).Value = "水"
$ws.Cells.Item(3049, 3).Value = 5
$ws.Cells.Item(3049, 4).Value = 201
$ws.Cells.Item(3050, 1).Value = "2026/12/30"
$ws.Cells.Item(3050, 2).Value = "水"
$ws.Cells.Item(3050, 3).Value = 8
$ws.Cells.Item(3050, 4).Value = 201
$ws.Cells.Item(3051, 1).Value = "2026/12/30"
$ws.Cells.Item(3051, 2).Value = "水"
$ws.Cells.Item(3051, 3).Value = 13
$ws.Cells.Item(3051, 4).Value = 201
$ws.Cells.Item(3052, 1).Value = "2026/12/30"
$ws.Cells.Item(3052, 2).Value = "水"
$ws.Cells.Item(3052, 3).Value = 16
$ws.Cells.Item(3052, 4).Value = 201
$ws.Cells.Item(3053, 1).Value = "2026/12/30"
$ws.Cells.Item(3053, 2).Value = "水"
$ws.Cells.Item(3053, 3).Value = 22
$ws.Cells.Item(3053, 4).Value = 201
$ws.Cells.Item(3054, 1).Value = "2026/12/31"
$ws.Cells.Item(3054, 2).Value = "木"
$ws.Cells.Item(3054, 3).Value = 2
$ws.Cells.Item(3054, 4).Value = 201
$ws.Cells.Item(3055, 1).Value = "2026/12/31"
$ws.Cells.Item(3055, 2).Value = "木"
$ws.Cells.Item(3055, 3).Value = 6
$ws.Cells.Item(3055, 4).Value = 201
$ws.Cells.Item(3056, 1).Value = "2026/12/31"
$ws.Cells.Item(3056, 2).Value = "木"
$ws.Cells.Item(3056, 3).Value = 10
$ws.Cells.Item(3056, 4).Value = 201
$ws.Cells.Item(3057, 1).Value = "2026/12/31"
$ws.Cells.Item(3057, 2).Value = "木"
$ws.Cells.Item(3057, 3).Value = 12
$ws.Cells.Item(3057, 4).Value = 201
$ws.Cells.Item(3058, 1).Value = "2026/12/31"
$ws.Cells.Item(3058, 2).Value = "木"
$ws.Cells.Item(3058, 3).Value = 14
$ws.Cells.Item(3058, 4).Value = 201
$ws.Cells.Item(3059, 1).Value = "2026/12/31"
$ws.Cells.Item(3059, 2).Value = "木"
$ws.Cells.Item(3059, 3).Value = 22
$ws.Cells.Item(3059, 4).Value = 201
$ws.Cells.Item(3060, 1).Value = "2027/01/01"
$ws.Cells.Item(3060, 2).Value = "金"
$ws.Cells.Item(3060, 3).Value = 2
$ws.Cells.Item(3060, 4).Value = 201
$ws.Cells.Item(3061, 1).Value = "2027/01/01"
$ws.Cells.Item(3061, 2).Value = "金"
$ws.Cells.Item(3061, 3).Value = 5
$ws.Cells.Item(3061, 4).Value = 201
$ws.Cells.Item(3062, 1).Value = "2027/01/01"
$ws.Cells.Item(3062, 2).Value = "金"
$ws.Cells.Item(3062, 3).Value = 13
$ws.Cells.Item(3062, 4).Value = 201
$ws.Cells.Item(3063, 1).Value = "2027/01/01"
$ws.Cells.Item(3063, 2).Value = "金"
$ws.Cells.Item(3063, 3).Value = 16
$ws.Cells.Item(3063, 4).Value = 201
$ws.Cells.Item(3064, 1).Value = "2027/01/01"
$ws.Cells.Item(3064, 2).Value = "金"
$ws.Cells.Item(3064, 3).Value = 19
$ws.Cells.Item(3064, 4).Value = 201
$ws.Cells.Item(3065, 1).Value = "2027/01/02"
$ws.Cells.Item(3065, 2).Value = "土"
$ws.Cells.Item(3065, 3).Value = 1
$ws.Cells.Item(3065, 4).Value = 201
$ws.Cells.Item(3066, 1).Value = "2027/01/02"
$ws.Cells.Item(3066, 2).Value = "土"
$ws.Cells.Item(3066, 3).Value = 5
$ws.Cells.Item(3066, 4).Value = 201
$ws.Cells.Item(3067, 1).Value = "2027/01/02"
$ws.Cells.Item(3067, 2).Value = "土"
$ws.Cells.Item(3067, 3).Value = 8
$ws.Cells.Item(3067, 4).Value = 201
$ws.Cells.Item(3068, 1).Value = "2027/01/02"
$ws.Cells.Item(3068, 2).Value = "土"
$ws.Cells.Item(3068, 3).Value = 13
$ws.Cells.Item(3068, 4).Value = 201
$ws.Cells.Item(3069, 1).Value = "2027/01/02"
$ws.Cells.Item(3069, 2).Value = "土"
$ws.Cells.Item(3069, 3).Value = 16
$ws.Cells.Item(3069, 4).Value = 201
$ws.Cells.Item(3070, 1).Value = "2027/01/02"
$ws.Cells.Item(3070, 2).Value = "土"
$ws.Cells.Item(3070, 3).Value = 19
$ws.Cells.Item(3070, 4).Value = 201
$ws.Cells.Item(3071, 1).Value = "2027/01/02"
$ws.Cells.Item(3071, 2).Value = "土"
$ws.Cells.Item(3071, 3).Value = 22
$ws.Cells.Item(3071, 4).Value = 201
$ws.Cells.Item(3072, 1).Value = "2027/01/03"
$ws.Cells.Item(3072, 2).Value = "日"
$ws.Cells.Item(3072, 3).Value = 1
$ws.Cells.Item(3072, 4).Value = 201
$ws.Cells.Item(3073, 1).Value = "2027/01/03"
$ws.Cells.Item(3073, 2).Value = "日"
$ws.Cells.Item(3073, 3).Value = 4
$ws.Cells.Item(3073, 4).Value = 201
$ws.Cells.Item(3074, 1).Value = "2027/01/03"
$ws.Cells.Item(3074, 2).Value = "日"
$ws.Cells.Item(3074, 3).Value = 7
$ws.Cells.Item(3074, 4).Value = 201
$ws.Cells.Item(3075, 1).Value = "2027/01/03"
$ws.Cells.Item(3075, 2).Value = "日"
$ws.Cells.Item(3075, 3).Value = 13
$ws.Cells.Item(3075, 4).Value = 201
$ws.Cells.Item(3076, 1).Value = "2027/01/03"
$ws.Cells.Item(3076, 2).Value = "日"
$ws.Cells.Item(3076, 3).Value = 16
$ws.Cells.Item(3076, 4).Value = 201
$ws.Cells.Item(3077, 1).Value = "2027/01/03"
$ws.Cells.Item(3077, 2).Value = "日"
$ws.Cells.Item(3077, 3).Value = 19
$ws.Cells.Item(3077, 4).Value = 201
$ws.Cells.Item(3078, 1).Value = "2027/01/03"
$ws.Cells.Item(3078, 2).Value = "日"
$ws.Cells.Item(3078, 3).Value = 22
$ws.Cells.Item(3078, 4).Value = 201
$ws.Cells.Item(3079, 1).Value = "2027/01/04"
$ws.Cells.Item(3079, 2).Value = "月"
$ws.Cells.Item(3079, 3).Value = 2
$ws.Cells.Item(3079, 4).Value = 201
$ws.Cells.Item(3080, 1).Value = "2027/01/04"
$ws.Cells.Item(3080, 2).Value = "月"
$ws.Cells.Item(3080, 3).Value = 4
$ws.Cells.Item(3080, 4).Value = 201
$ws.Cells.Item(3081, 1).Value = "2027/01/04"
$ws.Cells.Item(3081, 2).Value = "月"
$ws.Cells.Item(3081, 3).Value = 7
$ws.Cells.Item(3081, 4).Value = 201
$ws.Cells.Item(3082, 1).Value = "2027/01/04"
$ws.Cells.Item(3082, 2).Value = "月"
$ws.Cells.Item(3082, 3).Value = 13
$ws.Cells.Item(3082, 4).Value = 201
$ws.Cells.Item(3083, 1).Value = "2027/01/04"
$ws.Cells.Item(3083, 2).Value = "月"
$ws.Cells.Item(3083, 3).Value = 22
$ws.Cells.Item(3083, 4).Value = 201
$ws.Cells.Item(3084, 1).Value = "2027/01/05"
$ws.Cells.Item(3084, 2).Value = "火"
$ws.Cells.Item(3084, 3).Value = 2
$ws.Cells.Item(3084, 4).Value = 201
$ws.Cells.Item(3085, 1).Value = "2027/01/05"
$ws.Cells.Item(3085, 2).Value = "火"
$ws.Cells.Item(3085, 3).Value = 7
$ws.Cells.Item(3085, 4).Value = 201

$dateRange.Style = "Normal"
